$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Coverity finding was replaced with a different function/file pair.
$ws.Range("A6").Value = "flashImage"
$ws.Range("C6").Value = "/home/rdkv-core/cov/cov-analysis-linux64-2023.6.0/bin/device/entservices-softwareupdate/FirmwareUpdate/FirmwareUpdateImplementation.cpp"

# Updated line number for the finding in row 6.
$ws.Range("D6").Value = 351

# Move the active selection, mirroring where the editor last clicked.
$ws.Range("C11").Select()
